# Hortaliza, Macroferia Regional de Talca - Pepino ensalada
# Insert one new weekly data row at row 648 (pushing the existing rows
# 648-706 down to 649-707), then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 648 - shifts rows 648:706 down to 649:707,
# which grows the used range from A1:R706 to A1:R707.
$ws.Rows("648:648").Insert()

# Fill in the newly inserted row 648 with its data.
$ws.Cells.Item(648, 1).Value  = 5
$ws.Cells.Item(648, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(648, 3).Value  = "Maule"
$ws.Cells.Item(648, 4).Value  = 45132
$ws.Cells.Item(648, 5).Value  = 7
$ws.Cells.Item(648, 6).Value  = 100112043
$ws.Cells.Item(648, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(648, 8).Value  = "Sin especificar"
$ws.Cells.Item(648, 9).Value  = "Primera"
$ws.Cells.Item(648, 10).Value = 500
$ws.Cells.Item(648, 11).Value = 8000
$ws.Cells.Item(648, 12).Value = 8000
$ws.Cells.Item(648, 13).Value = 8000
$ws.Cells.Item(648, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(648, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(648, 16).Value = 133
$ws.Cells.Item(648, 17).Value = 60
$ws.Cells.Item(648, 18).Value = "Hortaliza"
